$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 185; this shifts the existing rows 185-209
# down to 186-210 and extends the sheet dimension to A1:T210.
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row 185 with the new weekly record
# (same "fixed" columns as the surrounding Membrillo / Vega Modelo de
# Temuco rows, new date + price data).
$ws.Range("A185").Value = 10
$ws.Range("B185").Value = "Vega Modelo de Temuco"
$ws.Range("C185").Value = "La Araucanía"
$ws.Range("D185").Value = 44776
$ws.Range("E185").Value = 9
$ws.Range("F185").Value = "Fruta"
$ws.Range("G185").Value = 100104
$ws.Range("H185").Value = "Frutos de pepita"
$ws.Range("I185").Value = 100104003
$ws.Range("J185").Value = "Membrillo"
$ws.Range("K185").Value = "Champion"
$ws.Range("L185").Value = "Primera"
$ws.Range("M185").Value = 80
$ws.Range("N185").Value = 10000
$ws.Range("O185").Value = 10000
$ws.Range("P185").Value = 10000
$ws.Range("Q185").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R185").Value = "Región de O'Higgins"
$ws.Range("S185").Value = 556
$ws.Range("T185").Value = 18
